# ---------------------------------------------------------------------------
# Daily Data sheet: append rows 346-462 (2023-12-11 .. 2024-04-05)
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$wsDaily = $wb.Worksheets.Item("Daily Data")

$dailyRows = @(
    @(346,45271,$null,$null),
    @(347,45272,$null,$null),
    @(348,45273,$null,$null),
    @(349,45274,$null,$null),
    @(350,45275,$null,$null),
    @(351,45276,$null,$null),
    @(352,45277,$null,$null),
    @(353,45278,$null,$null),
    @(354,45279,$null,$null),
    @(355,45280,$null,$null),
    @(356,45281,$null,$null),
    @(357,45282,$null,$null),
    @(358,45283,$null,$null),
    @(359,45284,$null,$null),
    @(360,45285,$null,$null),
    @(361,45286,$null,$null),
    @(362,45287,$null,$null),
    @(363,45288,$null,$null),
    @(364,45289,$null,$null),
    @(365,45290,$null,$null),
    @(366,45291,$null,$null),
    @(367,45292,0,0),
    @(368,45293,100,100),
    @(369,45294,100,100),
    @(370,45295,100,100),
    @(371,45296,100,100),
    @(372,45297,100,100),
    @(373,45298,0,0),
    @(374,45299,100,100),
    @(375,45300,100,100),
    @(376,45301,100,100),
    @(377,45302,100,100),
    @(378,45303,100,100),
    @(379,45304,100,100),
    @(380,45305,0,0),
    @(381,45306,100,100),
    @(382,45307,100,100),
    @(383,45308,100,100),
    @(384,45309,100,100),
    @(385,45310,100,100),
    @(386,45311,100,100),
    @(387,45312,0,0),
    @(388,45313,100,100),
    @(389,45314,100,100),
    @(390,45315,100,100),
    @(391,45316,100,100),
    @(392,45317,0,0),
    @(393,45318,100,100),
    @(394,45319,100,100),
    @(395,45320,100,100),
    @(396,45321,100,100),
    @(397,45322,100,100),
    @(398,45323,100,100),
    @(399,45324,100,100),
    @(400,45325,100,100),
    @(401,45326,0,0),
    @(402,45327,100,100),
    @(403,45328,100,100),
    @(404,45329,100,100),
    @(405,45330,100,100),
    @(406,45331,100,100),
    @(407,45332,100,100),
    @(408,45333,0,0),
    @(409,45334,100,100),
    @(410,45335,100,100),
    @(411,45336,100,100),
    @(412,45337,100,100),
    @(413,45338,100,100),
    @(414,45339,100,100),
    @(415,45340,0,0),
    @(416,45341,100,100),
    @(417,45342,100,100),
    @(418,45343,100,100),
    @(419,45344,100,100),
    @(420,45345,100,100),
    @(421,45346,100,100),
    @(422,45347,0,0),
    @(423,45348,100,100),
    @(424,45349,100,100),
    @(425,45350,100,100),
    @(426,45351,100,100),
    @(427,45352,100,100),
    @(428,45353,100,100),
    @(429,45354,0,0),
    @(430,45355,100,100),
    @(431,45356,100,100),
    @(432,45357,100,100),
    @(433,45358,100,100),
    @(434,45359,100,100),
    @(435,45360,100,100),
    @(436,45361,0,0),
    @(437,45362,100,100),
    @(438,45363,100,100),
    @(439,45364,100,100),
    @(440,45365,100,100),
    @(441,45366,100,100),
    @(442,45367,100,100),
    @(443,45368,0,0),
    @(444,45369,100,100),
    @(445,45370,100,100),
    @(446,45371,100,100),
    @(447,45372,100,100),
    @(448,45373,100,100),
    @(449,45374,100,100),
    @(450,45375,0,0),
    @(451,45376,0,0),
    @(452,45377,100,100),
    @(453,45378,100,100),
    @(454,45379,100,100),
    @(455,45380,100,100),
    @(456,45381,100,100),
    @(457,45382,0,0),
    @(458,45383,100,100),
    @(459,45384,100,100),
    @(460,45385,100,100),
    @(461,45386,100,100),
    @(462,45387,100,100)
)

$wsDaily.Range("A346:A462").NumberFormat = "yyyy/mm/dd"
foreach ($row in $dailyRows) {
    $r = $row[0]
    $wsDaily.Cells.Item($r, 1).Value = $row[1]
    if ($row[2] -ne $null) {
        $wsDaily.Cells.Item($r, 2).Value = $row[2]
        $wsDaily.Cells.Item($r, 3).Value = $row[3]
    }
}
# ---------------------------------------------------------------------------
# Weekly Data sheet: backfill C41:D45, fix C46:D49, append weeks for 2024
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Data")

$wsWeekly.Range("C41:D49").Value = 100

$weeklyMonths = @{21="2024-01"; 22="2024-02"; 23="2024-03"; 24="2024-04"}
$weeklyRows = @(
    @(50,21,"W1"),
    @(51,21,"W2"),
    @(52,21,"W3"),
    @(53,21,"W4"),
    @(54,22,"W1"),
    @(55,22,"W2"),
    @(56,22,"W3"),
    @(57,22,"W4"),
    @(58,23,"W1"),
    @(59,23,"W2"),
    @(60,23,"W3"),
    @(61,23,"W4"),
    @(62,24,"W1")
)

foreach ($row in $weeklyRows) {
    $r = $row[0]
    $monthKey = $row[1]
    $weekLabel = $row[2]
    $wsWeekly.Cells.Item($r, 1).Value = $weeklyMonths[$monthKey]
    $wsWeekly.Cells.Item($r, 2).Value = $weekLabel
    $wsWeekly.Cells.Item($r, 3).Value = 100
    $wsWeekly.Cells.Item($r, 4).Value = 100
}
